# Turn the "Input" sheet into a minimal, single-row RPA input template:
# keep the header row and exactly one example data row (F / Zalau / Vopsit),
# dropping the remaining sample rows entirely.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Input")

# Drop the extra sample rows (3-11) first so the shared-string table sheds
# every value that is no longer referenced anywhere in the workbook.
$ws.Range("A3:C11").ClearContents()

# Overwrite the remaining example row with the new sample values.
$ws.Range("A2").Value = "F"
$ws.Range("B2").Value = "Zalau"
$ws.Range("C2").Value = "Vopsit"

# Match the recorded selection left on the Input sheet.
$null = $ws.Range("C4").Select()
